$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "61.420.42"
$ws.Range("D3").Value = "2.885.60"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.25"
$ws.Range("E5").Value = "  -4.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.77"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "2.884.44"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.00"
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.01"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").Value = "3.362.99"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "61.416.65"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").Value = "2.878.65"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.45"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.83"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.28"
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.93"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("E26").Value = "  -10.13%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -5.41%  "
$ws.Range("E29").Value = "  +3.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("E32").Value = "  -7.55%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.42"
$ws.Range("E35").Value = "  -3.36%  "
$ws.Range("E36").Value = "  -3.17%  "
$ws.Range("E37").Value = "  -3.28%  "
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("E39").Value = "  -7.82%  "
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.20"
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("E42").Value = "  -2.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.26"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("E44").Value = "  -4.36%  "
$ws.Range("D45").Value = "2.708.84"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.08"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "343.10"
$ws.Range("E48").Value = "  -4.66%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("E51").Value = "  -4.43%  "
$ws.Range("ZZ1").Copy()
$ws.Range("D4,D5,D6,D8,D10,D14,D20,D21,D23,D24,D25,D30,D35,D41,D43,D46,D48").PasteSpecial(-4122)
$excel.CutCopyMode = $false
